# ajout de pikachu et mewtwo dans les herbes
$wb = $excel.ActiveWorkbook

$wsPokemon = $wb.Worksheets.Item("pokemon")
$wsTerrain = $wb.Worksheets.Item("element_terrain")

# update image paths on the "pokemon" sheet
$wsPokemon.Range("B2").Value = "assets/for_alex/pikachu.png"
$wsPokemon.Range("C2").Value = "assets/for_alex/mewtwo.png"

# update rarity (%) values, keeping them as plain text (no quote-prefix style)
$wsPokemon.Range("B3").NumberFormat = "@"
$wsPokemon.Range("B3").Value = "80.0"
$wsPokemon.Range("B3").Style = "Normal"

$wsPokemon.Range("C3").NumberFormat = "@"
$wsPokemon.Range("C3").Value = "20.0"
$wsPokemon.Range("C3").Style = "Normal"

# widen the new columns to fit the longer paths
$wsPokemon.Columns.Item(2).ColumnWidth = 32.166666666666664
$wsPokemon.Columns.Item(3).ColumnWidth = 35

# update selection on element_terrain sheet (no longer the active tab)
[void]$wsTerrain.Range("B8").Select()

# make "pokemon" the active sheet/tab and set its selection
[void]$wsPokemon.Activate()
[void]$wsPokemon.Range("C3").Select()
